$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two "MON Nov 20" / " 09:38:57 PST 2017" runs into one run
#    by replacing the combined text with itself (Find/Replace coalesces
#    the matched text into a single run).
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Execute("MON Nov 20 09:38:57 PST 2017", $true, $false, $false, $false, $false, $true, 1, $false, "MON Nov 20 09:38:57 PST 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Append the new "SAT Dec 16" purchase-details block right after the
#    paragraph that holds "Amount balance ... - 28357.0".
# ---------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$found = $findRng.Find.Execute("- 28357.0", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $findRng.Paragraphs(1)
$insertPoint = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)

$block = "`rSAT Dec 16 10:16:00 PST 2017`rPerson Name`t`t`t`t- HG`rBill number`t`t`t`t- 1970`r---------------------------------------------------------------`rItem Name`t`t`t`t- CARROT`rNumber of Pockets`t`t`t- 1`rNumber of KGs`t`t`t- 101`rRate`t`t`t`t`t- 38`rTotal Price`t`t`t`t- 3838.0`rAmount Received`t`t`t- 4000`rAmount balance`t`t`t- 28195.0`rAmount Received mode`t`t- CASH`r`r"

$insertPoint.InsertAfter($block)

Write-Output "inserted block"
